$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 361.75
$ws.Range("I12").Value = 419
$ws.Range("K12").Value = 419
$ws.Range("M12").Value = -249
$ws.Range("H13").Value = 4750
$ws.Range("J13").Value = 4750
$ws.Range("L13").Value = 4750
$ws.Range("N13").Value = -5088
$ws.Range("H19").Value = 1795
$ws.Range("I19").Value = 1795
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1795
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -1620
$ws.Range("N19").Value = ""
$ws.Range("H29").Value = 7333.3335
$ws.Range("I29").Value = 2000
$ws.Range("K29").Value = 6000
$ws.Range("M29").Value = -5719
$ws.Range("H33").Value = 277.84616
$ws.Range("I33").Value = 94.454544
$ws.Range("J33").Value = 1286.5
$ws.Range("K33").Value = 94.454544
$ws.Range("L33").Value = 1286.5
$ws.Range("M33").Value = 134.545456
$ws.Range("N33").Value = -1744.5
$ws.Range("H76").Value = 3750.9583
$ws.Range("I76").Value = 2999.889
$ws.Range("K76").Value = 2999.889
$ws.Range("M76").Value = -2684.889
$ws.Range("H79").Value = 3750.9583
$ws.Range("I79").Value = 2999.889
$ws.Range("K79").Value = 2999.889
$ws.Range("M79").Value = -1907.889
$ws.Range("H112").Value = 4218.375
$ws.Range("J112").Value = 4249.5713
$ws.Range("L112").Value = 12748.7139
$ws.Range("N112").Value = -14964.7139
$ws.Range("H116").Value = 46309
$ws.Range("J116").Value = 83333.336
$ws.Range("L116").Value = 83333.336
$ws.Range("N116").Value = -90217.336
$ws.Range("H137").Value = 21744392
$ws.Range("I137").Value = 62501190
$ws.Range("J137").Value = 7429.8
$ws.Range("K137").Value = 187503570
$ws.Range("L137").Value = 22289.4
$ws.Range("M137").Value = -187501020
$ws.Range("N137").Value = -27389.4
$ws.Range("H138").Value = 2780.2927
$ws.Range("I138").Value = 1616.9412
$ws.Range("K138").Value = 4850.8236
$ws.Range("M138").Value = 289.1764000000003
$ws.Range("H141").Value = 4901.926
$ws.Range("I141").Value = 1667.6
$ws.Range("K141").Value = 5002.799999999999
$ws.Range("M141").Value = 177.2000000000007

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 402.1111
$ws.Range("I5").Value = 231.28572
$ws.Range("K5").Value = 231.28572
$ws.Range("M5").Value = -119.28572
$ws.Range("H32").Value = 181622.88
$ws.Range("I32").Value = 248222.61
$ws.Range("J32").Value = 21000
$ws.Range("K32").Value = 248222.61
$ws.Range("L32").Value = 21000
$ws.Range("M32").Value = -247935.61
$ws.Range("N32").Value = -21574
$ws.Range("H61").Value = 2384178.8
$ws.Range("I61").Value = 3181.4595
$ws.Range("J61").Value = 20003558
$ws.Range("K61").Value = 3181.4595
$ws.Range("L61").Value = 20003558
$ws.Range("M61").Value = -2969.4595
$ws.Range("N61").Value = -20003982
$ws.Range("H122").Value = 2077.1
$ws.Range("I122").Value = 1697.6666
$ws.Range("K122").Value = 5092.9998
$ws.Range("M122").Value = -2642.9998
$ws.Range("H132").Value = 631435.9399999999
$ws.Range("I132").Value = 764695.75
$ws.Range("K132").Value = 2294087.25
$ws.Range("M132").Value = -2291557.25
$ws.Range("H136").Value = 2384178.8
$ws.Range("I136").Value = 3181.4595
$ws.Range("J136").Value = 20003558
$ws.Range("K136").Value = 9544.378499999999
$ws.Range("L136").Value = 60010674
$ws.Range("M136").Value = -6994.378499999999
$ws.Range("N136").Value = -60015774
$ws.Range("H137").Value = 85613.36
$ws.Range("J137").Value = 85613.36
$ws.Range("L137").Value = 85613.36
$ws.Range("N137").Value = -95813.36
$ws.Range("H138").Value = 114133.86
$ws.Range("J138").Value = 114133.86
$ws.Range("L138").Value = 114133.86
$ws.Range("N138").Value = -124413.86
$ws.Range("H139").Value = 120712.875
$ws.Range("J139").Value = 120712.875
$ws.Range("L139").Value = 120712.875
$ws.Range("N139").Value = -130992.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 402.1111
$ws.Range("I4").Value = 231.28572
$ws.Range("K4").Value = 231.28572
$ws.Range("M4").Value = -116.28572
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("H94").Value = 7324.8696
$ws.Range("I94").Value = 8375.157999999999
$ws.Range("J94").Value = 2336
$ws.Range("K94").Value = 8375.157999999999
$ws.Range("L94").Value = 2336
$ws.Range("M94").Value = -7924.157999999999
$ws.Range("N94").Value = -3238
$ws.Range("H134").Value = 3711148.2
$ws.Range("I134").Value = 5290.2144
$ws.Range("J134").Value = 9814914
$ws.Range("K134").Value = 15870.6432
$ws.Range("L134").Value = 29444742
$ws.Range("M134").Value = -13335.6432
$ws.Range("N134").Value = -29449812
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4277079.5
$ws.Range("I31").Value = 4277079.5
$ws.Range("K31").Value = 4277079.5
$ws.Range("M31").Value = -4276784.5
$ws.Range("H34").Value = 4277079.5
$ws.Range("I34").Value = 4277079.5
$ws.Range("K34").Value = 4277079.5
$ws.Range("M34").Value = -4276877.5
$ws.Range("H62").Value = 4999
$ws.Range("J62").Value = 4999
$ws.Range("L62").Value = 4999
$ws.Range("N62").Value = -6247
$ws.Range("H65").Value = 4999
$ws.Range("J65").Value = 4999
$ws.Range("L65").Value = 24995
$ws.Range("N65").Value = -31235
$ws.Range("H134").Value = 1825.6046
$ws.Range("I134").Value = 1448.7693
$ws.Range("J134").Value = 5499.75
$ws.Range("K134").Value = 4346.3079
$ws.Range("L134").Value = 16499.25
$ws.Range("M134").Value = -1811.3079
$ws.Range("N134").Value = -21569.25
$ws.Range("H140").Value = 85645.88
$ws.Range("J140").Value = 85645.88
$ws.Range("L140").Value = 85645.88
$ws.Range("N140").Value = -96005.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 247.41667
$ws.Range("J92").Value = 247.77777
$ws.Range("L92").Value = 743.33331
$ws.Range("N92").Value = -3239.33331
$ws.Range("H97").Value = 592.875
$ws.Range("J97").Value = 592.875
$ws.Range("L97").Value = 1778.625
$ws.Range("N97").Value = -2770.625
$ws.Range("H101").Value = 7704538.5
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 7704538.5
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 23113615.5
$ws.Range("M101").Value = ""
$ws.Range("N101").Value = -23118483.5
$ws.Range("H102").Value = 4867.7144
$ws.Range("I102").Value = 4888.385
$ws.Range("K102").Value = 14665.155
$ws.Range("M102").Value = -12231.155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = ""
$ws.Range("H14").Value = 1507250
$ws.Range("I14").Value = 1507250
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1507250
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1507082
$ws.Range("N14").Value = ""
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = ""
$ws.Range("H70").Value = 40818.91
$ws.Range("I70").Value = 32714.143
$ws.Range("J70").Value = 55002.25
$ws.Range("K70").Value = 32714.143
$ws.Range("L70").Value = 55002.25
$ws.Range("M70").Value = -32444.143
$ws.Range("N70").Value = -55542.25
$ws.Range("H73").Value = 40818.91
$ws.Range("I73").Value = 32714.143
$ws.Range("J73").Value = 55002.25
$ws.Range("K73").Value = 32714.143
$ws.Range("L73").Value = 55002.25
$ws.Range("M73").Value = -31778.143
$ws.Range("N73").Value = -56874.25
$ws.Range("H80").Value = 1782.8334
$ws.Range("I80").Value = 1739.6
$ws.Range("J80").Value = 1999
$ws.Range("K80").Value = 1739.6
$ws.Range("L80").Value = 1999
$ws.Range("M80").Value = -741.5999999999999
$ws.Range("N80").Value = -3995
$ws.Range("H83").Value = 1782.8334
$ws.Range("I83").Value = 1739.6
$ws.Range("J83").Value = 1999
$ws.Range("K83").Value = 8698
$ws.Range("L83").Value = 9995
$ws.Range("M83").Value = -3706
$ws.Range("N83").Value = -19979
$ws.Range("H132").Value = 12141.074
$ws.Range("I132").Value = 9150.375
$ws.Range("J132").Value = 36066.668
$ws.Range("K132").Value = 27451.125
$ws.Range("L132").Value = 108200.004
$ws.Range("M132").Value = -24921.125
$ws.Range("N132").Value = -113260.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4489.222
$ws.Range("I7").Value = 4299.75
$ws.Range("K7").Value = 4299.75
$ws.Range("M7").Value = -4187.75
$ws.Range("H126").Value = 4489.222
$ws.Range("I126").Value = 4299.75
$ws.Range("K126").Value = 12899.25
$ws.Range("M126").Value = -10429.25
$ws.Range("H132").Value = 8349513
$ws.Range("J132").Value = 3434.25
$ws.Range("L132").Value = 10302.75
$ws.Range("N132").Value = -15362.75
$ws.Range("H135").Value = 62365
$ws.Range("J135").Value = 62365
$ws.Range("L135").Value = 62365
$ws.Range("N135").Value = -72505
$ws.Range("H137").Value = 99999
$ws.Range("J137").Value = 99999
$ws.Range("L137").Value = 99999
$ws.Range("N137").Value = -110199

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 92199.62
$ws.Range("I122").Value = 8790
$ws.Range("K122").Value = 26370
$ws.Range("M122").Value = -23920
